$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (42-49) appended to the egg notes table
$data = @(
    @(1875, "highlight"),
    @(2525, "some food"),
    @(2413, "some food"),
    @(2017, "some food"),
    @(1913, "some food"),
    @(1880, "highlight on side"),
    @(2528, "some food"),
    @(1983, "some food")
)

$startRow = 42
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Update the view to match the post-edit state (active cell moves to B50)
$ws.Range("B50").Select()
